# Updated Global Glider Cal and Ingest sheets:
#   - Changed Cal scattering angle to 140
#   - Changed angular resolution to 1.13

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Asset_Cal_Info")

# CC_scattering_angle value (row 2, column F)
$ws.Range("F2").Value = 140

# CC_angular_resolution value (row 4, column F)
$ws.Range("F4").Value = 1.13

# Leave the sheet active with the cursor where the author left it
$ws.Activate()
$ws.Range("E15").Select() | Out-Null
